$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry "dwg_index" row, appended right after the last existing
# data row (54) -> new row 55.
$srcRow = 54
$newRow = 55

# Columns A, B and E on this sheet are formatted as Text (same as the
# existing "Part Number"/"Drawing Date" columns) so the new drawing index
# and test values are stored as literal text, not re-interpreted as a
# number/date.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("B" + $newRow).NumberFormat = "@"
$ws.Range("E" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "0018643"
$ws.Range("B" + $newRow).Value = "test"

# Carry the remaining columns over from the prior row, unchanged.
$ws.Range("C" + $newRow).Value = $ws.Range("C" + $srcRow).Text
$ws.Range("D" + $newRow).Value = $ws.Range("D" + $srcRow).Text
$ws.Range("E" + $newRow).Value = $ws.Range("E" + $srcRow).Text
$ws.Range("F" + $newRow).Value = $ws.Range("F" + $srcRow).Text

# Move the active selection down to the new row, scrolling the viewport
# so row 55 is visible.
$excel.Goto($ws.Range("A61"), $true)
